$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from an existing header cell (H1) to the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Data values for columns I (I0) and J (IF), rows 2-25
$data = @(
    @(11, 11),
    @(7, 8),
    @(8, 9),
    @(6, 7),
    @(7, 7),
    @(3, 5),
    @(4, 5),
    @(8, 9),
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(6, 7),
    @(5, 6),
    @(7, 7),
    @(9, 9),
    @(7, 7),
    @(7, 8),
    @(1, 1),
    @(7, 8),
    @(7, 7),
    @(6, 6),
    @(6, 6),
    @(5, 5),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
